# "Add files via upload" — re-saved workbook with an updated USD Amount
# figure for the Deposit/Crypto/Roobic row (T2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 235377
